# Apply weekly fruit/vegetable price update for Frambuesa (Vega Modelo de Temuco)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44616
$ws.Range("N2").Value = 3200
$ws.Range("O2").Value = 3200
$ws.Range("P2").Value = 3200
$ws.Range("S2").Value = 3200

# Row 3
$ws.Range("D3").Value = 44176
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 3000

# Row 4
$ws.Range("D4").Value = 44999
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 2500
$ws.Range("O4").Value = 2500
$ws.Range("P4").Value = 2500
$ws.Range("R4").Value = "Región de La Araucanía"
$ws.Range("S4").Value = 2500

# Row 5
$ws.Range("D5").Value = 44567
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 2400
$ws.Range("O5").Value = 2400
$ws.Range("P5").Value = 2400
$ws.Range("S5").Value = 2400

# Row 6
$ws.Range("D6").Value = 44214
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 1800
$ws.Range("O6").Value = 1800
$ws.Range("P6").Value = 1800
$ws.Range("S6").Value = 1800

# Row 8
$ws.Range("D8").Value = 44998
$ws.Range("N8").Value = 2500
$ws.Range("O8").Value = 2500
$ws.Range("P8").Value = 2500
$ws.Range("S8").Value = 2500

# Row 9
$ws.Range("D9").Value = 44551
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 4500
$ws.Range("O9").Value = 4500
$ws.Range("P9").Value = 4500
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 4500

# Row 11
$ws.Range("D11").Value = 44215
$ws.Range("M11").Value = 65
$ws.Range("N11").Value = 2800
$ws.Range("O11").Value = 2800
$ws.Range("P11").Value = 2800
$ws.Range("R11").Value = "Región de La Araucanía"
$ws.Range("S11").Value = 2800

# Row 12
$ws.Range("D12").Value = 44323
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 3200
$ws.Range("O12").Value = 3200
$ws.Range("P12").Value = 3200
$ws.Range("S12").Value = 3200

# Row 13
$ws.Range("D13").Value = 44574
$ws.Range("N13").Value = 3000
$ws.Range("O13").Value = 3000
$ws.Range("P13").Value = 3000
$ws.Range("S13").Value = 3000
